# Update the "two-digit divided by one-digit" answer table with a new
# generated set of problems, as described in the commit:
# "Update master to output generated at c8c62b6"
#
# The document contains a single 5-column table whose data rows are
# rows 1, 5, 9, 13 and 17 (the rows in between are blank spacer rows).
# We address each cell by (row, column) rather than by searching for its
# old text, because several of the new values collide with old values
# used elsewhere in the table (e.g. "36÷8=4, 4" and "62÷4=15, 2" each
# appear as both an old value in one cell and a new value in another),
# which would make a simple text Find/Replace ambiguous or order-sensitive.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $text) {
    $cell = $table.Cell($row, $col)
    $cell.Range.Text = $text
}

# Row 1 (first data row)
Set-CellText $t 1 1 "62÷5=12, 2"
Set-CellText $t 1 2 "36÷8=4, 4"
Set-CellText $t 1 3 "64÷9=7, 1"
Set-CellText $t 1 4 "97÷8=12, 1"
Set-CellText $t 1 5 "66÷7=9, 3"

# Row 5 (second data row)
Set-CellText $t 5 1 "35÷9=3, 8"
Set-CellText $t 5 2 "38÷6=6, 2"
Set-CellText $t 5 3 "57÷5=11, 2"
Set-CellText $t 5 4 "42÷8=5, 2"
Set-CellText $t 5 5 "50÷4=12, 2"

# Row 9 (third data row)
Set-CellText $t 9 1 "83÷5=16, 3"
Set-CellText $t 9 2 "52÷7=7, 3"
Set-CellText $t 9 3 "26÷6=4, 2"
Set-CellText $t 9 4 "17÷3=5, 2"
Set-CellText $t 9 5 "77÷5=15, 2"

# Row 13 (fourth data row)
Set-CellText $t 13 1 "53÷6=8, 5"
Set-CellText $t 13 2 "33÷4=8, 1"
Set-CellText $t 13 3 "99÷5=19, 4"
Set-CellText $t 13 4 "62÷4=15, 2"
Set-CellText $t 13 5 "52÷4=13, 0"

# Row 17 (fifth data row) - column 3 ("80÷8=10, 0") is unchanged per the diff.
Set-CellText $t 17 1 "84÷3=28, 0"
Set-CellText $t 17 2 "64÷2=32, 0"
Set-CellText $t 17 4 "64÷8=8, 0"
Set-CellText $t 17 5 "50÷5=10, 0"
